$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$s = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '47.144.13'
$ws.Range('D2').Style = $s
$ws.Range('E2').Value = '  +5.60%  '
$s = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.345.71'
$ws.Range('D3').Style = $s
$ws.Range('E3').Value = '  +4.72%  '
$s = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = $s
$ws.Range('E4').Value = '  -0.81%  '
$s = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.72'
$ws.Range('D5').Style = $s
$ws.Range('E5').Value = '  +0.43%  '
$s = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.96'
$ws.Range('D6').Style = $s
$ws.Range('E6').Value = '  +5.23%  '
$s = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.578'
$ws.Range('D7').Style = $s
$ws.Range('E7').Value = '  +1.87%  '
$ws.Range('E8').Value = '  -0.63%  '
$s = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.539'
$ws.Range('D9').Style = $s
$ws.Range('E9').Value = '  +5.40%  '
$s = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.97'
$ws.Range('D10').Style = $s
$ws.Range('E10').Value = '  +4.12%  '
$s = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0811'
$ws.Range('D11').Style = $s
$ws.Range('E11').Value = '  +1.71%  '
$s = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.46'
$ws.Range('D12').Style = $s
$ws.Range('E12').Value = '  +4.67%  '
$ws.Range('E13').Value = '  -0.32%  '
$s = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.700.60'
$ws.Range('D14').Style = $s
$ws.Range('E14').Value = '  +4.60%  '
$s = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.341.89'
$ws.Range('D15').Style = $s
$ws.Range('E15').Value = '  +0.45%  '
$s = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.21'
$ws.Range('D16').Style = $s
$ws.Range('E16').Value = '  +5.69%  '
$s = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.833'
$ws.Range('D17').Style = $s
$ws.Range('E17').Value = '  +0.99%  '
$s = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '46.917.89'
$ws.Range('D18').Style = $s
$ws.Range('E18').Value = '  +5.24%  '
$s = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.67'
$ws.Range('D19').Style = $s
$ws.Range('E19').Value = '  +17.35%  '
$s = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0952'
$ws.Range('D20').Style = $s
$ws.Range('E20').Value = '  +2.42%  '
$s = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.20'
$ws.Range('D21').Style = $s
$ws.Range('E21').Value = '  +0.90%  '
$s = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.46'
$ws.Range('D22').Style = $s
$ws.Range('E22').Value = '  +3.53%  '
$s = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '248.26'
$ws.Range('D23').Style = $s
$ws.Range('E23').Value = '  +4.81%  '
$s = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.00'
$ws.Range('D24').Style = $s
$ws.Range('E24').Value = '  +2.54%  '
$ws.Range('E25').Value = '  +2.40%  '
$ws.Range('E26').Value = '  -0.35%  '
$ws.Range('E27').Value = '  +16.55%  '
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('E29').Value = '  +2.13%  '
$s = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.21'
$ws.Range('D30').Style = $s
$ws.Range('E30').Value = '  +1.71%  '
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('E32').Value = '  +5.97%  '
$s = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '148.74'
$ws.Range('D33').Style = $s
$s = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.62'
$ws.Range('D34').Style = $s
$ws.Range('E34').Value = '  +0.22%  '
$s = $ws.Range('D35').Style
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.20'
$ws.Range('D35').Style = $s
$ws.Range('E35').Value = '  +1.82%  '
$s = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.114'
$ws.Range('D36').Style = $s
$ws.Range('E36').Value = '  +5.20%  '
$ws.Range('E37').Value = '  +1.64%  '
$s = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.83'
$ws.Range('D38').Style = $s
$ws.Range('E38').Value = '  -0.97%  '
$s = $ws.Range('D39').Style
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.03'
$ws.Range('D39').Style = $s
$ws.Range('E39').Value = '  +7.22%  '
$ws.Range('E40').Value = '  +6.63%  '
$s = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.43'
$ws.Range('D41').Style = $s
$ws.Range('E41').Value = '  +3.27%  '
$s = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '14.06'
$ws.Range('D42').Style = $s
$ws.Range('E42').Value = '  -6.00%  '
$ws.Range('E43').Value = '  -0.95%  '
$s = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.99'
$ws.Range('D44').Style = $s
$ws.Range('E44').Value = '  +13.06%  '
$s = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.846.48'
$ws.Range('D45').Style = $s
$ws.Range('E45').Value = '  +2.12%  '
$s = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '87.31'
$ws.Range('D46').Style = $s
$ws.Range('E46').Value = '  +7.49%  '
$s = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '76.01'
$ws.Range('D47').Style = $s
$ws.Range('E47').Value = '  +10.96%  '
$ws.Range('E48').Value = '  +6.12%  '
$s = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '99.44'
$ws.Range('D49').Style = $s
$ws.Range('E49').Value = '  +2.05%  '
$s = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.89'
$ws.Range('D50').Style = $s
$ws.Range('E50').Value = '  +1.65%  '
$s = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '55.50'
$ws.Range('D51').Style = $s
$ws.Range('E51').Value = '  +3.73%  '
